$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MuniEntryPleas")
$ws.Rows("9:10").Delete()
